$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = "Indisponibilidade de conexão com a Internet no campus"
$ws.Range("B4").Value = "Baixa"
$ws.Range("C4").Value = "Alto"
$ws.Range("D4").Value = "Comunicar a universidade sobre o ocorrido e utilizar solução secundária "

$ws.Range("A5").Select()
